$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '47.267.14'
$ws.Range('E2').Value = '  +5.83%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.512.83'
$ws.Range('E3').Value = '  +3.78%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').Value = '  -0.01%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '324.06'
$ws.Range('E5').Value = '  +2.22%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '105.51'
$ws.Range('E6').Value = '  +4.20%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.523'
$ws.Range('E7').Value = '  +2.15%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.999'
$ws.Range('E8').Value = '  -0.05%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.540'
$ws.Range('E9').Value = '  +2.58%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.70'
$ws.Range('E10').Value = '  +4.06%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0818'
$ws.Range('E11').Value = '  +2.48%  '

# Row 12
$ws.Range('E12').Value = '  +0.83%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '18.39'
$ws.Range('E13').Value = '  -1.18%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.17'
$ws.Range('E14').Value = '  +4.14%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.900.68'
$ws.Range('E15').Value = '  +3.69%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.511.94'
$ws.Range('E16').Value = '  +3.75%  '

# Row 17
$ws.Range('E17').Value = '  +2.17%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '47.190.00'
$ws.Range('E18').Value = '  +6.04%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.77'
$ws.Range('E19').Value = '  +4.67%  '

# Row 20
$ws.Range('E20').Value = '  +2.79%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0942'
$ws.Range('E21').Value = '  +2.79%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '70.90'
$ws.Range('E22').Value = '  +3.42%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '252.61'
$ws.Range('E23').Value = '  +4.45%  '

# Row 24
$ws.Range('E24').Value = '  +5.35%  '

# Row 25
$ws.Range('E25').Value = '  +2.99%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '26.40'
$ws.Range('E26').Value = '  +4.68%  '

# Row 27
$ws.Range('E27').Value = '  -0.03%  '

# Row 28
$ws.Range('B28').Value = 'Toncoin'
$ws.Range('C28').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.30'
$ws.Range('E28').Value = '  +0.58%  '

# Row 29
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.97'
$ws.Range('E29').Value = '  +5.03%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.10'
$ws.Range('E30').Value = '  +5.18%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.136'
$ws.Range('E31').Value = '  +8.59%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '49.75'
$ws.Range('E32').Value = '  +3.08%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '19.87'
$ws.Range('E33').Value = '  +0.90%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.32'
$ws.Range('E34').Value = '  +3.23%  '

# Row 35
$ws.Range('E35').Value = '  +2.02%  '

# Row 36
$ws.Range('E36').Value = '  +0.17%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.96'
$ws.Range('E37').Value = '  +4.11%  '

# Row 38
$ws.Range('E38').Value = '  +4.32%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.98'
$ws.Range('E39').Value = '  +4.56%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '123.73'
$ws.Range('E40').Value = '  -1.65%  '

# Row 41
$ws.Range('E41').Value = '  +2.24%  '

# Row 42
$ws.Range('E42').Value = '  +2.72%  '

# Row 43
$ws.Range('E43').Value = '  +2.54%  '

# Row 44
$ws.Range('E44').Value = '  +3.26%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.978.91'
$ws.Range('E45').Value = '  +2.03%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.02'
$ws.Range('E46').Value = '  +3.13%  '

# Row 47
$ws.Range('E47').Value = '  +0.55%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.80'
$ws.Range('E48').Value = '  +3.29%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.10'
$ws.Range('E49').Value = '  -0.27%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '5.39'
$ws.Range('E50').Value = '  +17.44%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '80.13'
$ws.Range('E51').Value = '  +6.49%  '
